# "Generate Report for Handback" - update localization-status report now that
# ea18cc3d-a63e-4ec2-8fc9-a03f72b78516.md has been handed back.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: status columns for the ea18cc3d row (row 3) ---
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: ea18cc3d row (row 3) ---
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K3").Value = "2016-08-18 02:42:26"
$wsZhCn.Range("P3").Value = ""

# --- de-de sheet: ea18cc3d row (row 3) ---
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K3").Value = "2016-08-18 02:42:37"
$wsDeDe.Range("P3").Value = ""

# The "Error Detail" column (P) no longer needs to be wide enough to hold the
# long stale-handback warning message, so shrink it back down on both
# language sheets.
$wsZhCn.Columns.Item(16).ColumnWidth = 12.8
$wsDeDe.Columns.Item(16).ColumnWidth = 12.8
